# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps on the Overview, zh-cn and de-de
# sheets to reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 4 is the c9c6d314-... entry, column G holds
# "Latest HO Xliff Generate Date".
$overview.Range("G4").Value = "2016-09-03 04:49:33"

# zh-cn sheet: row 4 is the c9c6d314-... entry.
#   H4 = Correspond Handoff Datetime
#   K4 = Correspond Handback DateTime
$zhcn.Range("H4").Value = "2016-09-03 04:49:28"
$zhcn.Range("K4").Value = "2016-09-03 04:49:47"

# de-de sheet: row 4 is the c9c6d314-... entry.
#   H4 = Correspond Handoff Datetime
#   K4 = Correspond Handback DateTime
$dede.Range("H4").Value = "2016-09-03 04:49:33"
$dede.Range("K4").Value = "2016-09-03 04:49:54"
